$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C2:C527) from 45181 to 45182
$ws.Range("C2:C527").Value = 45182
